# "updated enrollment confirmation layout"
#
# Mark a batch of checklist rows as completed: stamp the "X" / date
# columns (A:C) that already mark other rows as done, onto rows
# 14,15,16,17,18,21,22,23,24,29,40 - and extend the yellow-highlight
# formatting from column D out to column E on rows 17/18 where the
# checklist text lives in E instead of D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$completedDate = 45868   # 2025-07-30 serial date, same as other "done" rows

# Row 9 is a fully-styled template for the "completed" look:
#   A = "X" (yellow fill), B = date (yellow fill), C = blank (yellow, bold),
#   D = label cell (yellow fill)
$templateRange = "A9:D9"

$targetRows = 14,15,16,17,18,21,22,23,24,29,40

foreach ($r in $targetRows) {
    $ws.Range($templateRange).Copy() | Out-Null
    $ws.Range("A$r`:D$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("A$r").Value = "X"
    $ws.Range("B$r").Value = $completedDate
}

# Rows 17 & 18 carry their checklist text in column E (not D), so extend
# the same yellow styling from D onto E for those two rows.
$ws.Range("D9").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null

# Row 23's C cell ends up styled like the A/D cells (plain yellow) rather
# than the bold "CC" style used elsewhere, so fix it up to match.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Move the active selection to F9, matching the saved view state.
$ws.Range("F9").Select()
